$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the PCR master-mix quantities (N= cycles, gDNA amount, total rxn volume)
$ws.Range("C2").Value = 33
$ws.Range("C8").Value = 5
$ws.Range("C12").Value = 100

# Update the two small side-note rows (index table just below the PCR calc table)
$ws.Range("G14").Value = 90
$ws.Range("B15").Value = 18
$ws.Range("E15").Value = "F primer + gDNA"
$ws.Range("G15").Value = 10

# Move the active selection
$ws.Range("B16").Select()
